$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 1
    9  = 2
    10 = -7
    11 = -4
    13 = 0
    14 = -3
    17 = 2
    19 = 0
    27 = -2
    30 = 1
    31 = -2
    32 = 0
    45 = -2
    51 = -5
    52 = -4
    53 = -3
    54 = -1
    60 = 4
    63 = 0
    66 = -2
    69 = 3
    75 = -1
    77 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
